$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.217.00"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "1.655.21"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  -0.62%  "
$ws.Range("D5").Value = "'219.11"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").Value = "'0.5229"
$ws.Range("E6").Value = "  -2.39%  "
$ws.Range("E7").Value = "  -0.55%  "
$ws.Range("D8").Value = "'0.2643"
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("D9").Value = "'0.06316"
$ws.Range("E9").Value = "  -1.29%  "
$ws.Range("D10").Value = "'20.63"
$ws.Range("E10").Value = "  -1.79%  "
$ws.Range("D11").Value = "'0.07767"
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").Value = "'4.571"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "1.674.41"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("D14").Value = "1.883.96"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").Value = "'0.5627"
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("D16").Value = "0.0₅8078"
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("D17").Value = "'65.28"
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("D18").Value = "26.217.30"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("D20").Value = "'4.721"
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("D21").Value = "'192.82"
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("D22").Value = "'10.26"
$ws.Range("D23").Value = "'6.028"
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("D24").Value = "'1.005"
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("D25").Value = "'144.45"
$ws.Range("E25").Value = "  -0.93%  "
$ws.Range("D26").Value = "'0.1202"
$ws.Range("E26").Value = "  -2.09%  "
$ws.Range("D27").Value = "'7.246"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").Value = "'15.99"
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("D29").Value = "'1.504"
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("D30").Value = "'0.05594"
$ws.Range("E30").Value = "  -5.02%  "
$ws.Range("D31").Value = "'1.278"
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("D32").Value = "'3.491"
$ws.Range("E32").Value = "  -2.61%  "
$ws.Range("D33").Value = "'3.379"
$ws.Range("E33").Value = "  +2.28%  "
$ws.Range("E34").Value = "  -1.38%  "
$ws.Range("D35").Value = "'2.801"
$ws.Range("E35").Value = "  -1.39%  "
$ws.Range("D36").Value = "'0.9450"
$ws.Range("E36").Value = "  -2.72%  "
$ws.Range("D37").Value = "'2.402"
$ws.Range("E37").Value = "  -1.22%  "
$ws.Range("D38").Value = "'0.5743"
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").Value = "'5.960"
$ws.Range("E40").Value = "  +1.41%  "
$ws.Range("D41").Value = "'2.591"
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("D42").Value = "'0.8471"
$ws.Range("E42").Value = "  -2.24%  "
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("D44").Value = "1.037.07"
$ws.Range("E44").Value = "  -3.65%  "
$ws.Range("D45").Value = "'102.47"
$ws.Range("E45").Value = "  -1.83%  "
$ws.Range("D46").Value = "1.794.59"
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("D47").Value = "'58.30"
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("D48").Value = "0.0₈105"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("D49").Value = "'0.9999"
$ws.Range("E49").Value = "  -1.67%  "
$ws.Range("D50").Value = "'0.05318"
$ws.Range("E50").Value = "  +2.89%  "
$ws.Range("D51").Value = "'8.058"
$ws.Range("E51").Value = "  -0.16%  "

# Reset style on text-forced numeric-looking cells so no stray
# number-format / quote-prefix style sticks to them.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
